$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the 4 new columns -------------------------------------------------
# "Division", "Section Code" and "Section Name" are inserted right after the
# existing "Org" column (old column B becomes column E, etc.), and a "Status"
# column is inserted right after "Phase Name" (which, after the first insert,
# now lives in column H, so the new column lands at column I).
$ws.Columns("B:D").Insert()
$ws.Columns("I").Insert()

# Inserting whole columns also stamps inherited formatting onto every row in
# that column range, including the report-header rows above the table (rows
# 1-6), which only ever contained data in column A. Strip those accidental
# blank cells back out so rows 1-6 are untouched, matching the template.
$ws.Range("B1:D6").Clear()

# --- Fill in the new header cells on row 8 ------------------------------------
$ws.Range("B8").Value2 = "Division"
$ws.Range("C8").Value2 = "Section Code"
$ws.Range("D8").Value2 = "Section Name"
$ws.Range("I8").Value2 = "Status"

# --- Column widths for the newly inserted columns -----------------------------
# (New columns don't automatically pick up a width from their neighbours.)
$ws.Columns.Item(2).ColumnWidth = 23.2    # Division       -> ~24.03
$ws.Columns.Item(3).ColumnWidth = 15.33   # Section Code   -> ~16.26
$ws.Columns.Item(4).ColumnWidth = 15.33   # Section Name   -> ~16.26
$ws.Columns.Item(9).ColumnWidth = 13.83   # Status         -> ~14.72

# --- Clean up the stray formatted-but-empty cells that used to trail off to
# the right of the header row (AMG8:AMJ8 originally); after the 4 column
# inserts they now sit 4 columns further out, at AMK:AMN. Removing those
# columns collapses the sheet dimension back down to A1:U8.
$ws.Columns("AMK:AMN").Delete()

# --- Restore the view/selection state -----------------------------------------
$ws.Range("Q8").Select()
